$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A167").Value = "Plan Escala;  `$146.917;  IVA: `$27.914;  `$174.831;  25-03-21"
$ws.Range("A168").Value = "Propiedad destacada x30;  `$587.668;  IVA: `$111.657;  Total: `$699.325;  25-03-21"
$ws.Range("A169").Value = "Propiedad destacada x30;  `$587.706;  IVA: `$111.664;  Total: `$699.370;  26-03-21"
$ws.Range("A170").Value = "Plan Escala;  `$146.926;  IVA: `$27.916;  `$174.842;  26-03-21"
